$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the stray "_GoBack" bookmark from the subtitle paragraph
#    (it is re-added at the very end of the document instead, see
#    below). Word itself only ever keeps one "_GoBack" bookmark, so
#    this mirrors the last-edit marker moving to the new last edit.
# ------------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# ------------------------------------------------------------------
# 2. The document currently ends with a single empty paragraph right
#    before the sectPr. Insert all of the new closing content *before*
#    that paragraph (leaving it untouched/empty) so that it keeps
#    being a clean, run-less paragraph we can drop the _GoBack
#    bookmark into afterwards.
# ------------------------------------------------------------------
$tail = $d.Paragraphs.Item($d.Paragraphs.Count)
$tailRange = $tail.Range
$tailRange.InsertBefore("Hier onder kunt u live zien hoe de informatie die de Tracker doorstuurt wordt verwerkt:`rWat voor code is gebruikt?`rOm de code zo kort mogelijk en te houden terwijl de functionaliteit hier niet aan lijdt hebben wij de code gebaseerd rondom de …`rWie gebruikt de Customer Tracker?`rOmdat de Customer Tracker nog een redelijk nieuw soort product is wordt het nog niet veel gebruikt.`rOnder de Huidige gebruikers valt de Mediamarkt op Beurs.`rAls u zich hieraan wilt toevoegen kunt u hier contact met ons opnemen om een bestelling te plaatsen of om te onderhandelen over een contract.`r")

# ------------------------------------------------------------------
# 3. Apply the "Kop1" heading style to the two new headings, and
#    underline "hier" in the last new paragraph. Addressed by
#    counting back from the trailing (still untouched, still empty)
#    paragraph so this does not depend on fragile absolute offsets.
# ------------------------------------------------------------------
$total = $d.Paragraphs.Count

$headingCode = $d.Paragraphs.Item($total - 6)
$headingCode.Style = "Kop1"

$headingUsers = $d.Paragraphs.Item($total - 4)
$headingUsers.Style = "Kop1"

$alsPara = $d.Paragraphs.Item($total - 1)
$alsStart = $alsPara.Range.Start
$prefix = "Als u zich hieraan wilt toevoegen kunt u "
$hierStart = $alsStart + $prefix.Length
$hierEnd = $hierStart + "hier".Length
$hierRange = $d.Range($hierStart, $hierEnd)
$hierRange.Font.Underline = 1

# ------------------------------------------------------------------
# 4. Re-add the "_GoBack" bookmark, now around the trailing empty
#    paragraph (still untouched, still run-less).
# ------------------------------------------------------------------
$tail = $d.Paragraphs.Item($d.Paragraphs.Count)
$d.Bookmarks.Add("_GoBack", $tail.Range)
